# working on payment list
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add "gdid" column cells (reusing the existing shared string "gdid")
$ws.Range("G15").Value = "gdid"
$ws.Range("H19").Value = "gdid"
$ws.Range("H22").Value = "gdid"
$ws.Range("H24").Value = "gdid"
$ws.Range("I36").Value = "gdid"
$ws.Range("I44").Value = "gdid"
$ws.Range("G49").Value = "gdid"
$ws.Range("H53").Value = "gdid"

# Update the view: scroll back to top-left A1 and move selection to K11
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K11").Select()
